# Apply updated odds values to Sheet1, reflecting the 2024-12-05 FlashScore data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.4
$ws.Range("I3").Value = 2.7
$ws.Range("Z3").Value = 29
$ws.Range("AA3").Value = 19
$ws.Range("AW3").Value = 201
$ws.Range("BB3").Value = 41
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.98
$ws.Range("G5").Value = 2.1
$ws.Range("I5").Value = 3.6
$ws.Range("Q5").Value = 2.01
$ws.Range("R5").Value = 1.89
$ws.Range("U5").Value = 1.75
$ws.Range("V5").Value = 2
$ws.Range("X5").Value = 10
$ws.Range("AG5").Value = 201
$ws.Range("AJ5").Value = 12
$ws.Range("AY5").Value = 19
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 3
$ws.Range("G7").Value = 1.4
$ws.Range("H7").Value = 4.33
$ws.Range("I7").Value = 8.5
$ws.Range("J7").Value = 1.95
$ws.Range("U7").Value = 2.5
$ws.Range("V7").Value = 1.5
$ws.Range("X7").Value = 5.5
$ws.Range("Z7").Value = 8.5
$ws.Range("AB7").Value = 41
$ws.Range("AF7").Value = 101
$ws.Range("AH7").Value = 15
$ws.Range("AN7").Value = 3.1
$ws.Range("AX7").Value = 9
$ws.Range("G8").Value = 1.38
$ws.Range("H8").Value = 4.5
$ws.Range("I8").Value = 8.5
$ws.Range("J8").Value = 1.95
$ws.Range("L8").Value = 9
$ws.Range("AD8").Value = 9
$ws.Range("AF8").Value = 101
$ws.Range("AJ8").Value = 26
$ws.Range("AK8").Value = 101
$ws.Range("AL8").Value = 67
$ws.Range("AN8").Value = 3.1
$ws.Range("AV8").Value = 81
$ws.Range("BA8").Value = 251
$ws.Range("BB8").Value = 301
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("Q9").Value = 2.35
$ws.Range("R9").Value = 1.57
$ws.Range("Q10").Value = 2.6
$ws.Range("R10").Value = 1.48
$ws.Range("G12").Value = 5
$ws.Range("I12").Value = 1.62
$ws.Range("J12").Value = 5.5
$ws.Range("O12").Value = 1.25
$ws.Range("P12").Value = 3.75
$ws.Range("Q12").Value = 1.9
$ws.Range("R12").Value = 1.95
$ws.Range("U12").Value = 1.91
$ws.Range("V12").Value = 1.91
$ws.Range("Y12").Value = 17
$ws.Range("AE12").Value = 17
$ws.Range("AG12").Value = 301
$ws.Range("AH12").Value = 7
$ws.Range("AI12").Value = 7.5
$ws.Range("AK12").Value = 12
$ws.Range("AN12").Value = 7
$ws.Range("AO12").Value = 29
$ws.Range("AQ12").Value = 101
$ws.Range("AR12").Value = 126
$ws.Range("AS12").Value = 251
$ws.Range("AU12").Value = 8.5
$ws.Range("AX12").Value = 3.6
$ws.Range("BA12").Value = 26
$ws.Range("BC12").Value = 151
$ws.Range("M13").Value = 1.04
$ws.Range("N13").Value = 13
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 4
$ws.Range("Q13").Value = 1.75
$ws.Range("R13").Value = 2.05
$ws.Range("AC13").Value = 12
$ws.Range("AE13").Value = 21
$ws.Range("G14").Value = 5.75
$ws.Range("I14").Value = 1.55
$ws.Range("L14").Value = 2.1
$ws.Range("X14").Value = 29
$ws.Range("AC14").Value = 13
$ws.Range("AH14").Value = 8
$ws.Range("AI14").Value = 8
$ws.Range("AR14").Value = 101
$ws.Range("Q15").Value = 1.7
$ws.Range("R15").Value = 2.1
$ws.Range("G16").Value = 1.9
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 2.63
$ws.Range("L16").Value = 4.75
$ws.Range("N16").Value = 7.5
$ws.Range("W16").Value = 6.5
$ws.Range("X16").Value = 8.5
$ws.Range("Z16").Value = 17
$ws.Range("AD16").Value = 6
$ws.Range("AF16").Value = 51
$ws.Range("AH16").Value = 10
$ws.Range("AK16").Value = 41
$ws.Range("AO16").Value = 11
$ws.Range("AQ16").Value = 41
$ws.Range("AY16").Value = 23
$ws.Range("M17").Value = 1.02
$ws.Range("N17").Value = 19
$ws.Range("W17").Value = 8.5
$ws.Range("AB17").Value = 34
$ws.Range("AH17").Value = 29
$ws.Range("BB17").Value = 301
$ws.Range("G20").Value = 2.65
$ws.Range("H20").Value = 3
$ws.Range("I20").Value = 2.72
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 3.3
$ws.Range("M20").Value = 1.09
$ws.Range("N20").Value = 6
$ws.Range("O20").Value = 1.42
$ws.Range("P20").Value = 2.67
$ws.Range("Q20").Value = 2.22
$ws.Range("R20").Value = 1.6
$ws.Range("S20").Value = 1.45
$ws.Range("T20").Value = 2.55
$ws.Range("U20").Value = 1.88
$ws.Range("V20").Value = 1.82
$ws.Range("W20").Value = 7.2
$ws.Range("X20").Value = 12.5
$ws.Range("AB20").Value = 37
$ws.Range("AC20").Value = 6
$ws.Range("AE20").Value = 15
$ws.Range("AF20").Value = 80
$ws.Range("AG20").Value = 700
$ws.Range("AH20").Value = 7.4
$ws.Range("AJ20").Value = 10
$ws.Range("AL20").Value = 25
$ws.Range("AM20").Value = 37
$ws.Range("AN20").Value = 4.5
$ws.Range("AP20").Value = 23
$ws.Range("AS20").Value = 300
$ws.Range("AT20").Value = 2.55
$ws.Range("AU20").Value = 7
$ws.Range("AV20").Value = 65
$ws.Range("AX20").Value = 4.55
$ws.Range("AY20").Value = 15
$ws.Range("AZ20").Value = 23
$ws.Range("BA20").Value = 70
$ws.Range("BB20").Value = 110
$ws.Range("BC20").Value = 300
$ws.Range("M22").Value = 1.05
$ws.Range("N22").Value = 8.5
$ws.Range("Q22").Value = 1.9
$ws.Range("R22").Value = 1.9
$ws.Range("Q23").Value = 1.88
$ws.Range("R23").Value = 1.93
$ws.Range("G24").Value = 1.36
$ws.Range("I24").Value = 7
$ws.Range("J24").Value = 1.83
$ws.Range("Z24").Value = 9
$ws.Range("AD24").Value = 9.5
$ws.Range("AE24").Value = 21
$ws.Range("AI24").Value = 41
$ws.Range("AK24").Value = 81
$ws.Range("AO24").Value = 6.5
$ws.Range("AP24").Value = 19
$ws.Range("AQ24").Value = 17
$ws.Range("BA24").Value = 151
$ws.Range("G25").Value = 1.22
$ws.Range("H25").Value = 5.2
$ws.Range("I25").Value = 12.5
$ws.Range("J25").Value = 1.62
$ws.Range("L25").Value = 9.5
$ws.Range("N25").Value = 16.5
$ws.Range("Q25").Value = 1.55
$ws.Range("R25").Value = 2.15
$ws.Range("S25").Value = 1.28
$ws.Range("T25").Value = 3.45
$ws.Range("U25").Value = 2.07
$ws.Range("W25").Value = 7.1
$ws.Range("Z25").Value = 7
$ws.Range("AA25").Value = 10.75
$ws.Range("AB25").Value = 32
$ws.Range("AC25").Value = 13.5
$ws.Range("AD25").Value = 11
$ws.Range("AE25").Value = 26
$ws.Range("AH25").Value = 32
$ws.Range("AI25").Value = 110
$ws.Range("AJ25").Value = 40
$ws.Range("AK25").Value = 500
$ws.Range("AL25").Value = 175
$ws.Range("AM25").Value = 120
$ws.Range("AO25").Value = 5.2
$ws.Range("AQ25").Value = 12.5
$ws.Range("AU25").Value = 9.25
$ws.Range("AX25").Value = 11.5
$ws.Range("AY25").Value = 75
$ws.Range("AZ25").Value = 60
$ws.Range("BB25").Value = 500
